$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.196.69'
$ws.Range("E2").Value = '  +3.76%  '

$ws.Range("D3").Value = '1.739.97'
$ws.Range("E3").Value = '  +3.65%  '

$ws.Range("E4").Value = '  +0.96%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3812'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3562'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.70%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.51'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.214'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07629'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.467'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.06%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.09%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.101'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.85%  '

$ws.Range("D16").Value = '1.749.66'
$ws.Range("E16").Value = '  +4.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001148'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.004'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06745'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '85.81'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.487'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.65%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.48%  '

$ws.Range("D24").Value = '25.243.68'
$ws.Range("E24").Value = '  +4.10%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.456'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.851'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.20%  '

$ws.Range("D29").Value = '1.941.99'
$ws.Range("E29").Value = '  +4.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.196'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +20.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.108'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.231'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.18'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +13.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.804'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08825'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.732'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.81%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02488'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.90%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06669'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.89%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.278'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2249'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.281'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6548'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.54%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.23%  '

$ws.Range("E45").Value = '  +0.94%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6261'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.888'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.175'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.55%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '131.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.98%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07383'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.46%  '
